$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the full used range first so leftover cells (e.g. old column B rows 8-9,
# old column E/F rows 5-6) from the previous, longer lists don't linger.
$ws.Range("A1:I12").ClearContents()

# Column A - Morphology (abbreviations added)
$ws.Range("A1").Value = "Morphology"
$ws.Range("A2").Value = "Branched (Br)"
$ws.Range("A3").Value = "Cushion-like (Cushion)"
$ws.Range("A4").Value = "Digitate (Dig)"
$ws.Range("A5").Value = "Encrusting (Enc)"
$ws.Range("A6").Value = "Filamentous (Fil)"
$ws.Range("A7").Value = "Foliose (Fol)"
$ws.Range("A8").Value = "Massive (Mas)"
$ws.Range("A9").Value = "Mushroom (Mush)"
$ws.Range("A10").Value = "Polypoid (Poly)"
$ws.Range("A11").Value = "Spherical (Sph)"
$ws.Range("A12").Value = "Stolonial (Stol)"

# Column B - Phyla (replaces Taxonomic Group list)
$ws.Range("B1").Value = "Phyla"
$ws.Range("B2").Value = "Chlorophyta"
$ws.Range("B3").Value = "Cnidaria"
$ws.Range("B4").Value = "Cyanobacteria"
$ws.Range("B5").Value = "Phaeophyta"
$ws.Range("B6").Value = "Porifera"
$ws.Range("B7").Value = "Rhodophyta"

# Column C - Calcification (abbreviations added)
$ws.Range("C1").Value = "Calcification"
$ws.Range("C2").Value = "Non-calcified (NC)"
$ws.Range("C3").Value = "Articulated (AC)"
$ws.Range("C4").Value = "Non-articulated (Non-AC)"
$ws.Range("C5").Value = "Hermatypic (Herm)"

# Column D - Energetic Resource (moved here from old column H)
$ws.Range("D1").Value = "Energetic Resource"
$ws.Range("D2").Value = "Autotrophy"
$ws.Range("D3").Value = "Heterotrophy"
$ws.Range("D4").Value = "Mixotrophy"

# Column E - Symbiosis (moved here from old column D)
$ws.Range("E1").Value = "Symbiosis"
$ws.Range("E2").Value = "Asymbiotic"
$ws.Range("E3").Value = "Symbiodineaceae"
$ws.Range("E4").Value = "Cyanobacterium"

# Column F - Maximum Size (moved here from old column E)
$ws.Range("F1").Value = "Maximum Size"
$ws.Range("F2").Value = "<10 cm"
$ws.Range("F3").Value = "10-20 cm"
$ws.Range("F4").Value = "21-50 cm"
$ws.Range("F5").Value = "51-100 cm"
$ws.Range("F6").Value = "> 100 cm"

# Column G - Growth rate (moved here from old column F)
$ws.Range("G1").Value = "Growth rate"
$ws.Range("G2").Value = "Very slow (<1 cm/yr)"
$ws.Range("G3").Value = "Slow (1 cm/yr)"
$ws.Range("G4").Value = "Moderate (>1 cm/yr)"
$ws.Range("G5").Value = "Fast (5-10 cm/yr)"
$ws.Range("G6").Value = "Very fast (>10 cm/yr)"

# Column H - Feeding Mode (moved here from old column G)
$ws.Range("H1").Value = "Feeding Mode"
$ws.Range("H2").Value = "Photosynthesis"
$ws.Range("H3").Value = "Active filter-feeding (pumping)"
$ws.Range("H4").Value = "Passive filter-feeding"

# Column I - Life Span (unchanged)
$ws.Range("I1").Value = "Life Span"
$ws.Range("I2").Value = "Annual"
$ws.Range("I3").Value = "Perennial"

# Update the active selection to match the saved view state
$ws.Range("D3").Select()
